$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.802.49'
$ws.Range('E2').Value = '  +0.35%  '
$ws.Range('D3').Value = '2.564.49'
$ws.Range('E3').Value = '  +0.56%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '302.82'
$ws.Range('E5').Value = '  +2.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '97.10'
$ws.Range('E6').Value = '  +6.23%  '
$ws.Range('E7').Value = '  +0.47%  '
$ws.Range('E8').Value = '  +0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.547'
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('E11').Value = '  +0.57%  '
$ws.Range('E12').Value = '  +10.20%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '7.68'
$ws.Range('E13').Value = '  +0.78%  '
$ws.Range('D14').Value = '2.528.12'
$ws.Range('E14').Value = '  +0.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.881'
$ws.Range('E15').Value = '  +2.00%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.47'
$ws.Range('E16').Value = '  +2.77%  '
$ws.Range('D17').Value = '42.833.99'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '13.39'
$ws.Range('E18').Value = '  +6.83%  '
$ws.Range('D19').Value = '0.0₃0986'
$ws.Range('E19').Value = '  +2.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.63'
$ws.Range('E20').Value = '  +1.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.51'
$ws.Range('E21').Value = '  -0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '258.18'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  +2.49%  '
$ws.Range('E24').Value = '  -1.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '28.30'
$ws.Range('E25').Value = '  -4.00%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '39.38'
$ws.Range('E27').Value = '  +9.71%  '
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('E29').Value = '  -1.50%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.01'
$ws.Range('E30').Value = '  +2.14%  '
$ws.Range('E31').Value = '  +4.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.75'
$ws.Range('E32').Value = '  +1.49%  '
$ws.Range('E33').Value = '  -0.26%  '
$ws.Range('B34').Value = 'EnergySwap'
$ws.Range('C34').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '27.05'
$ws.Range('E34').Value = '  +11.92%  '
$ws.Range('B35').Value = 'Hedera'
$ws.Range('C35').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0801'
$ws.Range('E35').Value = '  +1.19%  '
$ws.Range('B36').Value = 'LidoDAOToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.31'
$ws.Range('E36').Value = '  -3.21%  '
$ws.Range('E37').Value = '  +0.93%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '17.99'
$ws.Range('E38').Value = '  +13.62%  '
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('E40').Value = '  +1.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '2.06'
$ws.Range('E41').Value = '  +30.72%  '
$ws.Range('B42').Value = 'NEARProtocol'
$ws.Range('C42').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.38'
$ws.Range('E42').Value = '  -1.29%  '
$ws.Range('B43').Value = 'VeChain'
$ws.Range('C43').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0307'
$ws.Range('E43').Value = '  -0.41%  '
$ws.Range('D44').Value = '2.062.55'
$ws.Range('E44').Value = '  -0.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.998'
$ws.Range('E45').Value = '  -0.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '88.71'
$ws.Range('E46').Value = '  +5.35%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.30'
$ws.Range('E47').Value = '  +6.16%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.65'
$ws.Range('E48').Value = '  +10.98%  '
$ws.Range('D49').Value = '2.812.93'
$ws.Range('E49').Value = '  +0.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '104.09'
$ws.Range('E50').Value = '  +0.96%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.190'
$ws.Range('E51').Value = '  +2.58%  '
